# Add Precision, Recall, F1-Score columns to the model-comparison sheet,
# inserted right after "Accuracy" and before the existing "AUC" column.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new columns where the old "AUC" column (C) used to be.
# This pushes AUC/TP/TN/FP/FN from C:G to F:J and keeps their values/styles.
$ws.Range("C1:E1").EntireColumn.Insert()

# New header labels
$ws.Range("C1").Value = "Precision"
$ws.Range("D1").Value = "Recall"
$ws.Range("E1").Value = "F1-Score"

# Decision Tree (row 2)
$ws.Range("C2").Value = 0.7659574468085106
$ws.Range("D2").Value = 0.72
$ws.Range("E2").Value = 0.7422680412371134

# Random Forest (row 3)
$ws.Range("C3").Value = 0.7924528301886793
$ws.Range("D3").Value = 0.84
$ws.Range("E3").Value = 0.8155339805825242

# Logistic Regression (row 4)
$ws.Range("C4").Value = 0.7547169811320755
$ws.Range("D4").Value = 0.8
$ws.Range("E4").Value = 0.7766990291262136
